$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.628.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.698.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5170"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.75%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06261"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07339"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.696.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.528"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5866"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.928.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008445"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -13.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.670.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.047"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.82%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "187.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.295"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.71%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.618"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1153"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.314"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05703"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.334"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.517"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.525"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.670"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.025"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6050"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.375"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.678"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.101.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01600"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8644"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.891"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.20%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.857.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.95%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000107"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.221"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.45%  "
